# Insert a new weekly price-report row for "Feria Lagunitas de Puerto Montt - Apio".
# The new record is inserted as row 258, pushing the existing rows 258:344 down
# to 259:345 (dimension grows from A1:R344 to A1:R345).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 258..344 down by one row.
$ws.Rows(258).Insert()

# Populate the newly inserted row 258 with the new observation.
$ws.Range("A258").Value = 4
$ws.Range("B258").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C258").Value = "Los Lagos"
$ws.Range("D258").Value = 44876
$ws.Range("E258").Value = 10
$ws.Range("F258").Value = 100112017
$ws.Range("G258").Value = "Apio"
$ws.Range("H258").Value = "Americana (o)"
$ws.Range("I258").Value = "Primera"
$ws.Range("J258").Value = 40
$ws.Range("K258").Value = 15000
$ws.Range("L258").Value = 15000
$ws.Range("M258").Value = 15000
$ws.Range("N258").Value = "$/docena de matas"
$ws.Range("O258").Value = "Región de Coquimbo"
$ws.Range("P258").Value = 2500
$ws.Range("Q258").Value = 6
$ws.Range("R258").Value = "Hortaliza"
